$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two data rows (rows 4 and 5) so the table shrinks to
# header + 2 rows (new dimension A1:G3).
$ws.Range("A4:G5").EntireRow.Delete()

# --- Row 2 -----------------------------------------------------------
# A2 holds a date-looking literal ("10.11.2022"); force Text format first
# so it is kept as a literal string instead of being parsed into a date
# serial number, then restore General so the cell keeps using the
# original "style 1" (numFmtId General) rather than minting a new style.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "10.11.2022"
$ws.Range("A2").NumberFormat = "General"

# B2 holds "25,36" which also looks numeric-ish; keep it literal text too.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "25,36"
$ws.Range("B2").NumberFormat = "General"

$ws.Range("C2").Value = "Общее"
# D2 is unchanged (already an empty value).
$ws.Range("E2").Value = "Ксеро"
$ws.Range("F2").Value = "Karpacz"
# G2 becomes blank (was "Бригада Миши").
$ws.Range("G2").ClearContents()

# --- Row 3 -----------------------------------------------------------
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "10.11.2022"
$ws.Range("A3").NumberFormat = "General"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "215"
$ws.Range("B3").NumberFormat = "General"

$ws.Range("C3").Value = "Люди"
$ws.Range("D3").Value = "Владислав"
$ws.Range("E3").Value = "Топливо"
$ws.Range("F3").Value = "Office"
# G3 becomes blank (was "Бригада Миши").
$ws.Range("G3").ClearContents()
